$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a value to be stored as literal Text (mirrors how this
# sheet was originally authored with inline/text cells, e.g. "110 €"
# instead of being auto-parsed as a currency number) and strip the
# quote-prefix style that gets attached so the cell stays style-less.
function Set-TextCell($row, $col, $text) {
    $c = $ws.Cells.Item($row, $col)
    $c.Value = "'" + $text
    $c.ClearFormats()
}

# Helper: an explicitly-present, but empty, Text cell (matches the
# <c t="inlineStr"/> placeholder cells used throughout this sheet).
function Set-EmptyText($row, $col) {
    $c = $ws.Cells.Item($row, $col)
    $c.Value = "'"
    $c.ClearFormats()
}

function Set-Number($row, $col, $num) {
    $ws.Cells.Item($row, $col).Value = $num
}

# Row 2
Set-Number 2 1 1
Set-TextCell 2 2 "110 €"
Set-TextCell 2 3 "110 €"
Set-EmptyText 2 4
Set-TextCell 2 5 "1/1/2025"

# Row 3
Set-Number 3 1 2
Set-TextCell 3 2 "60 €"
Set-TextCell 3 3 "60 €"
Set-EmptyText 3 4
Set-TextCell 3 5 "2/1/2025"

# Row 4
Set-Number 4 1 3
Set-TextCell 4 2 "100 €"
Set-EmptyText 4 3
Set-EmptyText 4 4
Set-TextCell 4 5 "6/1/2025"

# Row 5
Set-Number 5 1 4
Set-TextCell 5 2 "210 €"
Set-TextCell 5 3 "310 €"
Set-EmptyText 5 4
Set-TextCell 5 5 "6/1/2025"

# Row 6
Set-Number 6 1 5
Set-TextCell 6 2 "60 €"
Set-EmptyText 6 3
Set-EmptyText 6 4
Set-TextCell 6 5 "7/1/2025"

# Row 7
Set-Number 7 1 6
Set-TextCell 7 2 "50 €"
Set-TextCell 7 3 "110 €"
Set-EmptyText 7 4
Set-TextCell 7 5 "7/1/2025"

# Row 8
Set-Number 8 1 7
Set-TextCell 8 2 "80 €"
Set-EmptyText 8 3
Set-EmptyText 8 4
Set-TextCell 8 5 "8/1/2025"

# Row 9
Set-Number 9 1 8
Set-TextCell 9 2 "150 €"
Set-TextCell 9 3 "230 €"
Set-EmptyText 9 4
Set-TextCell 9 5 "8/1/2025"

# Row 10
Set-Number 10 1 9
Set-TextCell 10 2 "100 €"
Set-EmptyText 10 3
Set-EmptyText 10 4
Set-TextCell 10 5 "9/1/2025"

# Row 11
Set-Number 11 1 10
Set-TextCell 11 2 "240 €"
Set-TextCell 11 3 "340 €"
Set-EmptyText 11 4
Set-TextCell 11 5 "9/1/2025"

# Row 12
Set-Number 12 1 11
Set-TextCell 12 2 "300 €"
Set-EmptyText 12 3
Set-EmptyText 12 4
Set-TextCell 12 5 "10/1/2025"

# Row 13
Set-Number 13 1 12
Set-TextCell 13 2 "310 €"
Set-TextCell 13 3 "610 €"
Set-EmptyText 13 4
Set-TextCell 13 5 "10/1/2025"

# Row 14
Set-Number 14 1 13
Set-TextCell 14 2 "110 €"
Set-TextCell 14 3 "110 €"
Set-EmptyText 14 4
Set-TextCell 14 5 "14/1/2025"

# Row 15
Set-Number 15 1 14
Set-TextCell 15 2 "80 €"
Set-EmptyText 15 3
Set-EmptyText 15 4
Set-TextCell 15 5 "16/1/2025"

# Row 16
Set-Number 16 1 15
Set-TextCell 16 2 "80 €"
Set-TextCell 16 3 "160 €"
Set-EmptyText 16 4
Set-TextCell 16 5 "16/1/2025"

# Row 17
Set-Number 17 1 16
Set-TextCell 17 2 "100 €"
Set-TextCell 17 3 "100 €"
Set-EmptyText 17 4
Set-TextCell 17 5 "21/1/2025"

# Row 18
Set-Number 18 1 17
Set-TextCell 18 2 "90 €"
Set-TextCell 18 3 "90 €"
Set-EmptyText 18 4
Set-TextCell 18 5 "22/1/2025"

# Row 19 (new)
Set-Number 19 1 18
Set-TextCell 19 2 "70 €"
Set-TextCell 19 3 "70 €"
Set-EmptyText 19 4
Set-TextCell 19 5 "23/1/2025"

# Row 20 (new; overwrites the old "TOTAL MES" row)
Set-Number 20 1 19
Set-TextCell 20 2 "110 €"
Set-TextCell 20 3 "110 €"
Set-EmptyText 20 4
Set-TextCell 20 5 "24/1/2025"

# Row 21 (new)
Set-Number 21 1 20
Set-TextCell 21 2 "60 €"
Set-TextCell 21 3 "60 €"
Set-EmptyText 21 4
Set-TextCell 21 5 "27/1/2025"

# Row 22 (new)
Set-Number 22 1 21
Set-TextCell 22 2 "50 €"
Set-TextCell 22 3 "50 €"
Set-EmptyText 22 4
Set-TextCell 22 5 "28/1/2025"

# Row 23 (new)
Set-Number 23 1 22
Set-TextCell 23 2 "60 €"
Set-EmptyText 23 3
Set-EmptyText 23 4
Set-TextCell 23 5 "29/1/2025"

# Row 24 (new)
Set-Number 24 1 23
Set-TextCell 24 2 "50 €"
Set-TextCell 24 3 "110 €"
Set-EmptyText 24 4
Set-TextCell 24 5 "29/1/2025"

# Row 25: blank spacer row
Set-EmptyText 25 1

# Row 26: new TOTAL row
$ws.Cells.Item(26, 1).Value = "TOTAL Enero"
Set-EmptyText 26 2
Set-TextCell 26 3 "2630 €"
Set-EmptyText 26 4
Set-EmptyText 26 5

# Row 27: blank spacer row
Set-EmptyText 27 1
